$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 269.72726
$ws.Range("I2").Value = 269.72726
$ws.Range("K2").Value = 269.72726
$ws.Range("M2").Value = -156.72726
$ws.Range("H9").Value = 7302.7856
$ws.Range("I9").Value = 9226.362999999999
$ws.Range("K9").Value = 9226.362999999999
$ws.Range("M9").Value = -9057.362999999999
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").Value = ""
$ws.Range("H103").Value = 1497.4
$ws.Range("J103").Value = 749.5
$ws.Range("L103").Value = 2248.5
$ws.Range("N103").Value = -3420.5
$ws.Range("H125").Value = 8351000
$ws.Range("I125").Value = 1695407.2
$ws.Range("K125").Value = 15258664.8
$ws.Range("M125").Value = -15256204.8
$ws.Range("H127").Value = 1649.6666
$ws.Range("J127").Value = 966.3333
$ws.Range("L127").Value = 2898.9999
$ws.Range("N127").Value = -12818.9999
$ws.Range("H132").Value = 3539.037
$ws.Range("I132").Value = 3539.037
$ws.Range("K132").Value = 10617.111
$ws.Range("M132").Value = -8087.110999999999
$ws.Range("H137").Value = 7158.3335
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 7158.3335
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 21475.0005
$ws.Range("M137").Value = ""
$ws.Range("N137").Value = -26575.0005

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 33666.332
$ws.Range("I22").Value = 499.5
$ws.Range("J22").Value = 100000
$ws.Range("K22").Value = 499.5
$ws.Range("L22").Value = 100000
$ws.Range("M22").Value = -200.5
$ws.Range("N22").Value = -100598
$ws.Range("H32").Value = 24997.324
$ws.Range("I32").Value = 9891.833000000001
$ws.Range("J32").Value = 52884.383
$ws.Range("K32").Value = 9891.833000000001
$ws.Range("L32").Value = 52884.383
$ws.Range("M32").Value = -9604.833000000001
$ws.Range("N32").Value = -53458.383
$ws.Range("H45").Value = 3254
$ws.Range("I45").Value = 3254
$ws.Range("K45").Value = 3254
$ws.Range("M45").Value = -2877
$ws.Range("H63").Value = 4331.6665
$ws.Range("I63").Value = 4331.6665
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 4331.6665
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -3645.6665
$ws.Range("N63").Value = ""
$ws.Range("H66").Value = 4331.6665
$ws.Range("I66").Value = 4331.6665
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 21658.3325
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -18226.3325
$ws.Range("N66").Value = ""
$ws.Range("H97").Value = 1413.3889
$ws.Range("I97").Value = 1307.3077
$ws.Range("J97").Value = 1689.2
$ws.Range("K97").Value = 1307.3077
$ws.Range("L97").Value = 1689.2
$ws.Range("M97").Value = -811.3077000000001
$ws.Range("N97").Value = -2681.2
$ws.Range("H122").Value = 1071.738
$ws.Range("I122").Value = 1040.325
$ws.Range("K122").Value = 3120.975
$ws.Range("M122").Value = -670.9750000000004

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2834.6667
$ws.Range("I20").Value = 1744.2727
$ws.Range("K20").Value = 1744.2727
$ws.Range("M20").Value = -1497.2727
$ws.Range("H94").Value = 710.913
$ws.Range("I94").Value = 710.913
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 710.913
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -259.913
$ws.Range("N94").Value = ""

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 38200
$ws.Range("I25").Value = 17600
$ws.Range("J25").Value = 100000
$ws.Range("K25").Value = 17600
$ws.Range("L25").Value = 100000
$ws.Range("M25").Value = -17426
$ws.Range("N25").Value = -100348
$ws.Range("H50").Value = 50000
$ws.Range("I50").Value = 50000
$ws.Range("K50").Value = 50000
$ws.Range("M50").Value = -49375
$ws.Range("H51").Value = 50000
$ws.Range("J51").Value = 50000
$ws.Range("L51").Value = 50000
$ws.Range("N51").Value = -51472
$ws.Range("H59").Value = 99999.664
$ws.Range("J59").Value = 99999.664
$ws.Range("L59").Value = 99999.664
$ws.Range("N59").Value = -102289.664
$ws.Range("H60").Value = 40833.168
$ws.Range("J60").Value = 40833.168
$ws.Range("L60").Value = 40833.168
$ws.Range("N60").Value = -41855.168
$ws.Range("H61").Value = 50000
$ws.Range("J61").Value = 50000
$ws.Range("L61").Value = 50000
$ws.Range("N61").Value = -50696
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = ""
$ws.Range("H141").Value = 573303.8
$ws.Range("J141").Value = 617593.1
$ws.Range("L141").Value = 617593.1
$ws.Range("N141").Value = -627953.1

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 536977.5
$ws.Range("I4").Value = 652883.0600000001
$ws.Range("K4").Value = 1958649.18
$ws.Range("M4").Value = -1958537.18
$ws.Range("H55").Value = 4200
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").Value = ""
$ws.Range("H88").Value = 15249.25
$ws.Range("J88").Value = 15249.25
$ws.Range("L88").Value = 45747.75
$ws.Range("N88").Value = -46603.75
$ws.Range("H91").Value = 15249.25
$ws.Range("J91").Value = 15249.25
$ws.Range("L91").Value = 45747.75
$ws.Range("N91").Value = -48711.75
$ws.Range("H141").Value = 2790.4666
$ws.Range("I141").Value = 2497.2856
$ws.Range("J141").Value = 6895
$ws.Range("K141").Value = 7491.8568
$ws.Range("L141").Value = 20685
$ws.Range("M141").Value = -2311.8568
$ws.Range("N141").Value = -31045

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1424.375
$ws.Range("I102").Value = 1413.5714
$ws.Range("K102").Value = 1413.5714
$ws.Range("M102").Value = 208.4286
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = ""
$ws.Range("H132").Value = 8587.75
$ws.Range("I132").Value = 8525.714
$ws.Range("K132").Value = 25577.142
$ws.Range("M132").Value = -23047.142

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 77241
$ws.Range("J6").Value = 77241
$ws.Range("L6").Value = 77241
$ws.Range("N6").Value = -77465
$ws.Range("H16").Value = 2599.4211
$ws.Range("I16").Value = 580.9091
$ws.Range("K16").Value = 580.9091
$ws.Range("M16").Value = -410.9091
$ws.Range("H22").Value = 1649.7
$ws.Range("I22").Value = 1008.4783
$ws.Range("K22").Value = 1008.4783
$ws.Range("M22").Value = -713.4783
$ws.Range("H27").Value = 1649.7
$ws.Range("I27").Value = 1008.4783
$ws.Range("K27").Value = 1008.4783
$ws.Range("M27").Value = -901.4783
$ws.Range("H68").Value = 2458426
$ws.Range("I68").Value = 3270864.5
$ws.Range("K68").Value = 3270864.5
$ws.Range("M68").Value = -3270115.5
$ws.Range("H71").Value = 2458426
$ws.Range("I71").Value = 3270864.5
$ws.Range("K71").Value = 16354322.5
$ws.Range("M71").Value = -16350578.5
$ws.Range("H95").Value = 19874.25
$ws.Range("J95").Value = 19874.25
$ws.Range("L95").Value = 19874.25
$ws.Range("N95").Value = -25366.25
$ws.Range("H122").Value = 4903.7075
$ws.Range("I122").Value = 4446.839
$ws.Range("K122").Value = 13340.517
$ws.Range("M122").Value = -10890.517
$ws.Range("H136").Value = 3049.1177
$ws.Range("I136").Value = 3066.5625
$ws.Range("J136").Value = 2770
$ws.Range("K136").Value = 9199.6875
$ws.Range("L136").Value = 8310
$ws.Range("M136").Value = -6649.6875
$ws.Range("N136").Value = -13410

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 60001164
$ws.Range("J4").Value = 75001350
$ws.Range("L4").Value = 75001350
$ws.Range("N4").Value = -75001576
$ws.Range("H81").Value = 2666.5
$ws.Range("I81").Value = 2289
$ws.Range("J81").Value = 4554
$ws.Range("K81").Value = 4578
$ws.Range("L81").Value = 9108
$ws.Range("M81").Value = -3517
$ws.Range("N81").Value = -11230
$ws.Range("H84").Value = 2666.5
$ws.Range("I84").Value = 2289
$ws.Range("J84").Value = 4554
$ws.Range("K84").Value = 22890
$ws.Range("L84").Value = 45540
$ws.Range("M84").Value = -17586
$ws.Range("N84").Value = -56148
$ws.Range("H93").Value = 68694
$ws.Range("J93").Value = 68694
$ws.Range("L93").Value = 68694
$ws.Range("N93").Value = -73686
$ws.Range("H132").Value = 5802.9067
$ws.Range("I132").Value = 3765.7334
$ws.Range("K132").Value = 11297.2002
$ws.Range("M132").Value = -8767.200199999999
